$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C (rows 2 through 238) holds a "Förändrad" date that was bumped
# from serial 45190 (2023-09-21) to serial 45192 (2023-09-23) for every row.
$lastRow = 238
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 3).Value = 45192
}
